# Small bug fix on student controller:
# The exported "demo" sample file was missing its header row, so Excel was
# treating row 1 as the first data record. Insert a real header row above
# the existing data (shifting all 30 student rows down by one) and label
# each column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 1; this pushes every
# existing row (and its cell styles) down by one without altering the data.
$ws.Rows(1).Insert()

# Populate the new header row with the column titles.
$ws.Range("A1").Value = "Admission Number"
$ws.Range("B1").Value = "Symbol Number/Roll Number"
$ws.Range("C1").Value = "Student Name"
$ws.Range("D1").Value = "Date Of Birth(BS)"
$ws.Range("E1").Value = "Religion"
$ws.Range("F1").Value = "Mobile No"
$ws.Range("G1").Value = "Email"
$ws.Range("H1").Value = "Admission Date"
$ws.Range("I1").Value = "Blood Group (O+, A+, B+, AB+, O-, A-, B-, AB-)"
$ws.Range("J1").Value = "Gender(Male/Female)"

# The mailto hyperlink that used to live on G1 (the first data row) now
# needs to point at G2, where that data moved to. Drop the stale link and
# re-create it on the new location, then restore the Hyperlink cell style
# (Hyperlinks.Add() on this engine always assigns a fresh style, so force
# it back to the shared built-in "Hyperlink" style that the sheet already
# uses elsewhere).
$ws.Range("G1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:email@email.com")
$ws.Range("G2").Style = "Hyperlink"

# Match the author's final selection.
$ws.Range("B2").Select()
